# Outlook login cred update - 4 Dec 2024
#
# The "Password" cell (B2) is being changed from the old mailto-hyperlinked
# value "Yankee@123456" to the new value "Welcome241029". The hyperlink
# relationship on B2 (rId2 -> mailto:Yankee@123456) is left untouched, but
# since the cell's visible text no longer matches the hyperlink's original
# text, Excel records that original text in the hyperlink's display string.
# The active selection also moves from B2 to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPassword = "Yankee@123456"
$newPassword = "Welcome241029"

# Preserve the hyperlink's original display text (the old password) on the
# B2 hyperlink before the cell text is overwritten. Iterating the
# Hyperlinks collection (rather than indexing a single Item(n) and writing
# through it) is what makes this an in-place update of the existing
# hyperlink entry instead of inserting a duplicate one.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 2) {
        $h.TextToDisplay = $oldPassword
    }
}

# Update the password cell itself to the new value.
$ws.Range("B2").Value = $newPassword

# Move the selection to D8.
$ws.Range("D8").Select() | Out-Null
